$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transitions")
$ws.Activate() | Out-Null

# Append three new ArmedTest transition rows to the bottom of the table
# (Origin State | Destination State | Trigger Source | Trigger Condition)
$newRows = @(
    @("ArmedTest", "Fault",    "internal"),
    @("ArmedTest", "Disarmed", "command"),
    @("ArmedTest", "Disarmed", "Arduino")
)

$startRow = 59
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Scroll the view down to the newly added rows and select the last new cell
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("C61").Select() | Out-Null
